# c01s03 - Linux Setup
#
# The deck has three "section divider" title slides whose titles read
# "Section 1 - Mac OS X Setup", "Section 2 - Linux Setup" and
# "Section 3 - Windows Setup" (slides 7, 8 and 9). The section numbers
# are bumped by one (a new section was inserted ahead of them), so:
#   Section 1 -> Section 2   (slide 7, Mac OS X Setup)
#   Section 2 -> Section 3   (slide 8, Linux Setup)
#   Section 3 -> Section 4   (slide 9, Windows Setup)
#
# Each title run is "Section N " (note trailing space) followed by a
# separate run for the en-dash and the rest of the subtitle, so only the
# digit + trailing space (characters 9-10) needs to be retyped; that is
# exactly what PowerPoint does when you double-click the digit and type
# the replacement, which splits "Section N " into "Section " + "N ".

$p = $ppt.ActivePresentation

$sections = @(
    @{ Slide = 7; New = "2 " },
    @{ Slide = 8; New = "3 " },
    @{ Slide = 9; New = "4 " }
)

foreach ($entry in $sections) {
    $slide = $p.Slides.Item($entry.Slide)
    $title = $slide.Shapes.Item(1)
    $tr = $title.TextFrame.TextRange
    # Characters 9-10 are the section digit plus the trailing space,
    # e.g. "Section [1 ]- Mac OS X Setup".
    $digit = $tr.Characters(9, 2)
    $digit.Text = $entry.New
}
